$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 10-13 (no longer part of the data; target cluster "MuSCs" removed as an option
# and sender/target re-shuffled per updated TPM-based NATMI computation)
$ws.Range("A10:T13").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd34"
$ws.Range("C2").Value = "Sell"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 227.282303
$ws.Range("H2").Value = 681.846909
$ws.Range("I2").Value = 0.6094595465130797
$ws.Range("J2").Value = 0.6094595465130795
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1797713333333334
$ws.Range("N2").Value = 0.5393140000000001
$ws.Range("O2").Value = 0.188800001120238
$ws.Range("P2").Value = 0.188800001120238
$ws.Range("Q2").Value = 40.85884265338067
$ws.Range("R2").Value = 367.729583880426
$ws.Range("S2").Value = 0.1150659630644092
$ws.Range("T2").Value = 0.1150659630644091

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd34"
$ws.Range("C3").Value = "Sell"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 227.282303
$ws.Range("H3").Value = 681.846909
$ws.Range("I3").Value = 0.6094595465130797
$ws.Range("J3").Value = 0.6094595465130795
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7724073333333333
$ws.Range("N3").Value = 2.317222
$ws.Range("O3").Value = 0.811199998879762
$ws.Range("P3").Value = 0.811199998879762
$ws.Range("Q3").Value = 175.5545175740887
$ws.Range("R3").Value = 1579.990658166798
$ws.Range("S3").Value = 0.4943935834486705
$ws.Range("T3").Value = 0.4943935834486704

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cd34"
$ws.Range("C4").Value = "Sell"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 138.990916
$ws.Range("H4").Value = 416.972748
$ws.Range("I4").Value = 0.3727053955221387
$ws.Range("J4").Value = 0.3727053955221385
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1797713333333334
$ws.Range("N4").Value = 0.5393140000000001
$ws.Range("O4").Value = 0.188800001120238
$ws.Range("P4").Value = 0.188800001120238
$ws.Range("Q4").Value = 24.98658229054134
$ws.Range("R4").Value = 224.879240614872
$ws.Range("S4").Value = 0.07036677909209851
$ws.Range("T4").Value = 0.07036677909209849

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd34"
$ws.Range("C5").Value = "Sell"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 138.990916
$ws.Range("H5").Value = 416.972748
$ws.Range("I5").Value = 0.3727053955221387
$ws.Range("J5").Value = 0.3727053955221385
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7724073333333333
$ws.Range("N5").Value = 2.317222
$ws.Range("O5").Value = 0.811199998879762
$ws.Range("P5").Value = 0.811199998879762
$ws.Range("Q5").Value = 107.3576027851173
$ws.Range("R5").Value = 966.218425066056
$ws.Range("S5").Value = 0.3023386164300401
$ws.Range("T5").Value = 0.3023386164300401

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cd34"
$ws.Range("C6").Value = "Sell"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.377814666666666
$ws.Range("H6").Value = 19.133444
$ws.Range("I6").Value = 0.0171021675827138
$ws.Range("J6").Value = 0.01710216758271379
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1797713333333334
$ws.Range("N6").Value = 0.5393140000000001
$ws.Range("O6").Value = 0.188800001120238
$ws.Range("P6").Value = 0.188800001120238
$ws.Range("Q6").Value = 1.146548246379556
$ws.Range("R6").Value = 10.318934217416
$ws.Range("S6").Value = 0.003228889258774863
$ws.Range("T6").Value = 0.003228889258774861

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cd34"
$ws.Range("C7").Value = "Sell"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.377814666666666
$ws.Range("H7").Value = 19.133444
$ws.Range("I7").Value = 0.0171021675827138
$ws.Range("J7").Value = 0.01710216758271379
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7724073333333333
$ws.Range("N7").Value = 2.317222
$ws.Range("O7").Value = 0.811199998879762
$ws.Range("P7").Value = 0.811199998879762
$ws.Range("Q7").Value = 4.926270819174221
$ws.Range("R7").Value = 44.33643737256799
$ws.Range("S7").Value = 0.01387327832393894
$ws.Range("T7").Value = 0.01387327832393893

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cd34"
$ws.Range("C8").Value = "Sell"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2733126666666666
$ws.Range("H8").Value = 0.819938
$ws.Range("I8").Value = 0.0007328903820679218
$ws.Range("J8").Value = 0.0007328903820679217
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1797713333333334
$ws.Range("N8").Value = 0.5393140000000001
$ws.Range("O8").Value = 0.188800001120238
$ws.Range("P8").Value = 0.188800001120238
$ws.Range("Q8").Value = 0.04913378250355556
$ws.Range("R8").Value = 0.442204042532
$ws.Range("S8").Value = 0.0001383697049554353
$ws.Range("T8").Value = 0.0001383697049554352

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cd34"
$ws.Range("C9").Value = "Sell"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2733126666666666
$ws.Range("H9").Value = 0.819938
$ws.Range("I9").Value = 0.0007328903820679218
$ws.Range("J9").Value = 0.0007328903820679217
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7724073333333333
$ws.Range("N9").Value = 2.317222
$ws.Range("O9").Value = 0.811199998879762
$ws.Range("P9").Value = 0.811199998879762
$ws.Range("Q9").Value = 0.2111087080262222
$ws.Range("R9").Value = 1.899978372236
$ws.Range("S9").Value = 0.0005945206771124866
$ws.Range("T9").Value = 0.0005945206771124865
